# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (A) stores one of three emoji markers describing a
# clinical trial's result status. Those emoji don't render reliably
# everywhere, so replace them with safer text/emoji equivalents:
#   📕 (red book)    -> "-3"
#   📘 (blue book)   -> "⚠️"
#   📙 (orange book) -> "+3"
#
# Only the "statut" column (A, below the header row) is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstDataRow = $used.Row + 1
$lastDataRow = $used.Row + $used.Rows.Count - 1

$statutRange = $ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($lastDataRow, 1))

# "-3" and "+3" look like numbers, so any cell about to receive one of
# those values must be switched to Text format first - otherwise Excel
# would reinterpret it as a number and the leading "+" would be lost.
# "⚠️" never parses as a number, so those cells are left alone.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Text
    if (($current -eq "📕") -or ($current -eq "📙")) {
        $cell.NumberFormat = "@"
    }
}

$statutRange.Replace("📕", "-3", 1, 1, $false, $false, $false)
$statutRange.Replace("📘", "⚠️", 1, 1, $false, $false, $false)
$statutRange.Replace("📙", "+3", 1, 1, $false, $false, $false)
